$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to stay text-formatted so numeric-looking strings
# (prices like "43.332.26", "0.470", percentages, etc.) are not
# auto-converted/truncated by Excel when assigned via .Value
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.332.26"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "2.231.94"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "258.28"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "0.623"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("E7").Value = "  +4.39%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").Value = "42.93"
$ws.Range("E10").Value = "  +4.09%  "
$ws.Range("D11").Value = "0.0924"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "7.08"
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "2.566.20"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "14.61"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "2.223.76"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "0.793"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "43.278.62"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").Value = "71.31"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "6.03"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("E22").Value = "  +5.90%  "
$ws.Range("D23").Value = "231.31"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").Value = "9.31"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "41.92"
$ws.Range("E26").Value = "  +6.77%  "
$ws.Range("D27").Value = "10.82"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").Value = "173.37"
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "20.47"
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.0868"
$ws.Range("E33").Value = "  +9.18%  "
$ws.Range("D34").Value = "5.25"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").Value = "0.122"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").Value = "0.0369"
$ws.Range("E36").Value = "  +13.13%  "
$ws.Range("D37").Value = "4.45"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  -5.44%  "
$ws.Range("D39").Value = "13.23"
$ws.Range("E39").Value = "  +7.24%  "
$ws.Range("D40").Value = "2.89"
$ws.Range("E40").Value = "  +18.81%  "
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").Value = "61.52"
$ws.Range("E43").Value = "  +3.07%  "
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "8.61"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "103.55"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "0.470"
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("D48").Value = "0.0980"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").Value = "1.48"
$ws.Range("E51").Value = "  +23.77%  "

Write-Output "Applied crypto price/volume updates"
